$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert new row 23: "Alexa Fluor 488 hCD1d unloaded" / "hCD1d unloaded"
# (a new NIH Tetramer Core reagent, alongside the existing PBS-57 loaded one)
# ---------------------------------------------------------------------------
$ws.Rows("23:23").Insert()
$ws.Range("A200").Copy()
$ws.Range("A23").PasteSpecial(-4122)

$ws.Range("B23").Value = "Alexa Fluor 488 hCD1d unloaded"
$ws.Range("E23").Value = "NIH Tetramer Core"
$ws.Range("F23").Value = "NA"
$ws.Range("J23").Value = "hCD1d unloaded"
$ws.Range("L23").Value = "NA"
$ws.Range("M23").Value = "Alexa Fluor 488"

# ---------------------------------------------------------------------------
# Insert new row 32: "Alexa Fluor 647 hMR1 6-FP" / "hMR1 6-FP"
# (a new NIH Tetramer Core reagent, alongside the existing 5-OP-RU loaded one)
# ---------------------------------------------------------------------------
$ws.Rows("32:32").Insert()
$ws.Range("A200").Copy()
$ws.Range("A32").PasteSpecial(-4122)

$ws.Range("B32").Value = "Alexa Fluor 647 hMR1 6-FP"
$ws.Range("E32").Value = "NIH Tetramer Core"
$ws.Range("F32").Value = "NA"
$ws.Range("J32").Value = "hMR1 6-FP"
$ws.Range("L32").Value = "NA"
$ws.Range("M32").Value = "Alexa Fluor 647"

# ---------------------------------------------------------------------------
# Restore the view: scroll position and active-cell selection
# ---------------------------------------------------------------------------
$ws.Range("A10").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("H27").Select() | Out-Null
